# Updates to 杭州-漫展信息.xlsx matching the commit diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) -- sheet index 1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 121
$ws1.Range("F4").Value = 889
$ws1.Range("F5").Value = 1062
$ws1.Range("F7").Value = 14
$ws1.Range("F8").Value = 655
$ws1.Range("F9").Value = 11941
$ws1.Range("F14").Value = 248
$ws1.Range("F16").Value = 1220
$ws1.Range("F17").Value = 190
$ws1.Range("F19").Value = 757
$ws1.Range("F20").Value = 670

$ws1.Range("D21").Value = "5号大街297号 盛泰开元名都大酒店"
$ws1.Range("F21").Value = 287

$ws1.Range("F22").Value = 2917
$ws1.Range("F24").Value = 3853
$ws1.Range("F25").Value = 1088
$ws1.Range("F26").Value = 840
$ws1.Range("F28").Value = 27
$ws1.Range("F30").Value = 1023
$ws1.Range("F32").Value = 93
$ws1.Range("F37").Value = 4382
$ws1.Range("F39").Value = 4503

$ws1.Range("E40").Value = "2024.11.09 09:30-11.10 17:00"
$ws1.Range("F40").Value = 5517

$ws1.Range("F42").Value = 124
$ws1.Range("F49").Value = 121

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) -- sheet index 2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F3").Value = 4170
$ws2.Range("F5").Value = 97
$ws2.Range("F12").Value = 898
$ws2.Range("F13").Value = 6

# Insert a new row before the existing row 25 ("...维也纳皇家交响乐团..."),
# which pushes that row down to row 26.
$ws2.Rows("25:25").Insert()

# Re-apply the bold/bordered "index" cell format to the new A25 (copied from
# the row that is now A26) before writing its value.
$ws2.Range("A26").Copy()
$ws2.Range("A25").PasteSpecial(-4122)

$ws2.Range("A25").Value = 24
$ws2.Range("B25").Value = "'2025-01-01"
$ws2.Range("B25").Style = "Normal"
$ws2.Range("C25").Value = "杭州·【早鸟优惠】大型正版授权互动卡通儿童剧《海底小纵队之深海探秘》"
$ws2.Range("D25").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws2.Range("E25").Value = "2025.01.01 10:30-01.01 11:40"
$ws2.Range("F25").Value = 0
$ws2.Range("G25").Value = 40
$ws2.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=92951"
$ws2.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202409/oZlaKX931727335820196.jpeg"

# Renumber the pushed-down row.
$ws2.Range("A26").Value = 25

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) -- sheet index 4
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F6").Value = 121
$ws4.Range("F8").Value = 889
$ws4.Range("F9").Value = 1062
$ws4.Range("F11").Value = 14
$ws4.Range("F12").Value = 655
$ws4.Range("F13").Value = 11941
$ws4.Range("F16").Value = 1220
$ws4.Range("F17").Value = 190
$ws4.Range("F19").Value = 4170
$ws4.Range("F20").Value = 757
$ws4.Range("F21").Value = 670
$ws4.Range("F22").Value = 2917
$ws4.Range("F24").Value = 3853
$ws4.Range("F25").Value = 3853
$ws4.Range("F26").Value = 1088
$ws4.Range("F27").Value = 840
$ws4.Range("F30").Value = 27
$ws4.Range("F33").Value = 1023
$ws4.Range("F35").Value = 93
$ws4.Range("F39").Value = 6
$ws4.Range("F41").Value = 124
